# Excel To CSV 工具修改
#
# - Re-type the "B1/C1" (and the type-suffix on "E1") header labels on both
#   sheets so the ".string" suffix becomes ".String" (capitalised), which is
#   what produced the new/renumbered shared-string table in the target file.
# - Switch the active/selected sheet from "EveryMonth" to "EveryDay".
#
# NOTE: the order in which new text is written matters, because it controls
# the order new entries are appended to xl/sharedStrings.xml. The target
# file appends them as: name.String, desc.String, triggerCondition.String,
# triggerConditions.String - i.e. EveryMonth's "triggerCondition.String" is
# written before EveryDay's "triggerConditions.String".

$wb = $excel.ActiveWorkbook

$wsDay = $wb.Worksheets.Item("EveryDay")
$wsMonth = $wb.Worksheets.Item("EveryMonth")

# EveryDay header row (A1:E1): id.int | name.String | desc.String | redreshTime.float.array | triggerConditions.String
$wsDay.Range("B1").Value = "name.String"
$wsDay.Range("C1").Value = "desc.String"

# EveryMonth header row (A1:E1): id.int | name.String | desc.String | redreshTime.float | triggerCondition.String
$wsMonth.Range("B1").Value = "name.String"
$wsMonth.Range("C1").Value = "desc.String"
$wsMonth.Range("E1").Value = "triggerCondition.String"

# EveryDay's E1 last, so "triggerConditions.String" is the final new shared string
$wsDay.Range("E1").Value = "triggerConditions.String"

# Make "EveryDay" the active/selected sheet (previously it was "EveryMonth")
$wsDay.Activate()
